$d = $word.ActiveDocument

$replacements = @(
    @("541×6=", "677×6="),
    @("712×3=", "466×6="),
    @("919×9=", "895×4="),
    @("550×9=", "630×8="),
    @("408×3=", "492×6="),
    @("181×7=", "104×9="),
    @("892×3=", "988×4="),
    @("499×6=", "698×7="),
    @("688×2=", "150×4="),
    @("846×7=", "858×2="),
    @("907×5=", "186×8="),
    @("948×2=", "164×8="),
    @("439×7=", "816×3="),
    @("937×3=", "923×7="),
    @("322×4=", "379×7="),
    @("699×9=", "443×5="),
    @("795×4=", "532×2="),
    @("339×7=", "828×7="),
    @("283×8=", "874×8="),
    @("258×9=", "338×2="),
    @("760×9=", "664×4="),
    @("106×4=", "548×5="),
    @("649×6=", "678×4="),
    @("749×8=", "887×4="),
    @("977×2=", "915×9=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
